$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Testy wydajnościowe" header cell to the new, broader label.
$ws.Range("C17").Value = "Testy wydajnościowe i obciążeniowe"

# The old "additional machine requirements" note is no longer needed here.
$ws.Range("C19").Value = ""

# Row 14 now autosizes to a slightly shorter custom height.
$ws.Rows.Item(14).RowHeight = 36.75

# Column D (4) is narrower than before.
$ws.Columns.Item(4).ColumnWidth = 18.7

# Move/refresh the active selection to H22.
$ws.Range("H22").Select() | Out-Null
